# Quarterly indexing esoteric bug-fix operation
#
# The Qn quarterly-error columns were off by one: a newly computed "Q0"
# error needs to be inserted immediately after the row label (column B),
# pushing every existing Q-error one column to the right. Rows that were
# already fully populated (B:K, 10 quarters) lose their former last
# value (column K) since there is no column beyond K for it to land in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = -0.6603092772102132
    3  = -0.15162438770796
    4  = -0.2053460154962278
    5  = 0.6162032393936197
    6  = 1.652643173475852
    7  = 0.3110387314724781
    8  = 0.2388379152847414
    9  = 0.6508000635779043
    10 = 0.2387740594105157
    11 = 0.3465902496671606
    12 = 0.00230005330798793
    13 = -0.1902738424076751
    14 = -0.3325070745318338
    15 = 0.1656141382254278
    16 = -0.09587373626955231
}

foreach ($row in 2..16) {
    # Shift the existing values one column to the right, working from the
    # rightmost populated column back down to column C so nothing gets
    # clobbered before it is read.
    for ($col = 11; $col -ge 3; $col--) {
        $srcVal = $ws.Cells.Item($row, $col - 1).Value2
        if ($srcVal -eq $null) {
            $ws.Cells.Item($row, $col).ClearContents()
        } else {
            $ws.Cells.Item($row, $col).Value2 = $srcVal
        }
    }

    # Place the newly computed Q0 error in column B.
    $ws.Cells.Item($row, 2).Value2 = $newValues[$row]
}
